$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 15
$ws.Range("A3").Value = 5
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 20
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 10
$ws.Range("A10").Value = 15
$ws.Range("A11").Value = 15
$ws.Range("A12").Value = 3
$ws.Range("A13").Value = 3

$ws.Rows.Item(3).RowHeight = 12.85

$ws.Range("A14").Select()
